$wb = $excel.ActiveWorkbook

# The diff updates F2 (1566 -> 1567) and F3 (92 -> 93) on both the
# "展览" sheet and the "全部类型" sheet (the two sheets that carry this
# data table; the other two sheets only contain a header row).
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1567
    $ws.Range("F3").Value = 93
}
